$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1543.8889
$ws.Range("I2").Value = 1656.5714
$ws.Range("J2").Value = 1149.5
$ws.Range("K2").Value = 1656.5714
$ws.Range("L2").Value = 1149.5
$ws.Range("M2").Value = -1543.5714
$ws.Range("N2").Value = -1375.5
$ws.Range("H70").Value = 3345.5386
$ws.Range("I70").Value = 3159.2
$ws.Range("J70").Value = 3966.6667
$ws.Range("K70").Value = 9477.599999999999
$ws.Range("L70").Value = 11900.0001
$ws.Range("M70").Value = -9207.599999999999
$ws.Range("N70").Value = -12440.0001
$ws.Range("H73").Value = 3345.5386
$ws.Range("I73").Value = 3159.2
$ws.Range("J73").Value = 3966.6667
$ws.Range("K73").Value = 9477.599999999999
$ws.Range("L73").Value = 11900.0001
$ws.Range("M73").Value = -8541.599999999999
$ws.Range("N73").Value = -13772.0001
$ws.Range("H80").Value = 2711.0688
$ws.Range("I80").Value = 2200.3
$ws.Range("K80").Value = 6600.900000000001
$ws.Range("M80").Value = -5602.900000000001
$ws.Range("H83").Value = 2711.0688
$ws.Range("I83").Value = 2200.3
$ws.Range("K83").Value = 19802.7
$ws.Range("M83").Value = -14810.7
$ws.Range("H86").Value = 2783.5
$ws.Range("J86").Value = 2750
$ws.Range("L86").Value = 2750
$ws.Range("N86").Value = -4996
$ws.Range("H89").Value = 2783.5
$ws.Range("J89").Value = 2750
$ws.Range("L89").Value = 13750
$ws.Range("N89").Value = -24982
$ws.Range("H113").Value = 7217.25
$ws.Range("I113").Value = 6900
$ws.Range("K113").Value = 6900
$ws.Range("M113").Value = -3646
$ws.Range("H116").Value = 11090.6
$ws.Range("I116").Value = 10213.8
$ws.Range("J116").Value = 12844.2
$ws.Range("K116").Value = 10213.8
$ws.Range("L116").Value = 12844.2
$ws.Range("M116").Value = -6771.799999999999
$ws.Range("N116").Value = -19728.2
$ws.Range("H133").Value = 120000
$ws.Range("J133").Value = 120000
$ws.Range("L133").Value = 120000
$ws.Range("N133").Value = -130120
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 4916.6665
$ws.Range("J12").Value = 4916.6665
$ws.Range("L12").Value = 4916.6665
$ws.Range("N12").Value = -5262.6665
$ws.Range("H14").Value = 949
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 949
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 949
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -1299
$ws.Range("H21").Value = 1766
$ws.Range("I21").Value = 350
$ws.Range("J21").Value = 2474
$ws.Range("K21").Value = 350
$ws.Range("L21").Value = 2474
$ws.Range("M21").Value = 24
$ws.Range("N21").Value = -3222
$ws.Range("H22").Value = 1738.6666
$ws.Range("I22").Value = 2108
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 2108
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -1809
$ws.Range("N22").Value = -1598
$ws.Range("H25").Value = 3928.8
$ws.Range("J25").Value = 9500
$ws.Range("L25").Value = 9500
$ws.Range("N25").Value = -10304
$ws.Range("H26").Value = 3203.5
$ws.Range("I26").Value = 3203.5
$ws.Range("K26").Value = 3203.5
$ws.Range("M26").Value = -2873.5
$ws.Range("H29").Value = 5000
$ws.Range("J29").Value = 5000
$ws.Range("L29").Value = 5000
$ws.Range("N29").Value = -5616
$ws.Range("H30").Value = 40001
$ws.Range("J30").Value = 40001
$ws.Range("L30").Value = 40001
$ws.Range("N30").Value = -40301
$ws.Range("H32").Value = 5804.675
$ws.Range("I32").Value = 5870.054
$ws.Range("K32").Value = 5870.054
$ws.Range("M32").Value = -5583.054
$ws.Range("H36").Value = 15500
$ws.Range("I36").Value = 15500
$ws.Range("K36").Value = 15500
$ws.Range("M36").Value = -15154
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H28").Value = 59771
$ws.Range("J28").Value = 59771
$ws.Range("L28").Value = 59771
$ws.Range("N28").Value = -60359
$ws.Range("H86").Value = 13942.308
$ws.Range("I86").Value = 7403.6113
$ws.Range("J86").Value = 28654.375
$ws.Range("K86").Value = 7403.6113
$ws.Range("L86").Value = 28654.375
$ws.Range("M86").Value = -6280.6113
$ws.Range("N86").Value = -30900.375
$ws.Range("H89").Value = 13942.308
$ws.Range("I89").Value = 7403.6113
$ws.Range("J89").Value = 28654.375
$ws.Range("K89").Value = 37018.0565
$ws.Range("L89").Value = 143271.875
$ws.Range("M89").Value = -31402.0565
$ws.Range("N89").Value = -154503.875
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 366.66666
$ws.Range("I17").Value = 400
$ws.Range("J17").Value = 300
$ws.Range("K17").Value = 1200
$ws.Range("L17").Value = 900
$ws.Range("M17").Value = -1031
$ws.Range("N17").Value = -1238
$ws.Range("H44").Value = 8036.3335
$ws.Range("I44").Value = 8299.237999999999
$ws.Range("J44").Value = 7116.1665
$ws.Range("K44").Value = 24897.714
$ws.Range("L44").Value = 21348.4995
$ws.Range("M44").Value = -24499.714
$ws.Range("N44").Value = -22144.4995
$ws.Range("H64").Value = 14521.375
$ws.Range("I64").Value = 1792.75
$ws.Range("K64").Value = 5378.25
$ws.Range("M64").Value = -5108.25
$ws.Range("H67").Value = 14521.375
$ws.Range("I67").Value = 1792.75
$ws.Range("K67").Value = 5378.25
$ws.Range("M67").Value = -4442.25
$ws.Range("H97").Value = 954.8
$ws.Range("J97").Value = 206.66667
$ws.Range("L97").Value = 620.00001
$ws.Range("N97").Value = -1612.00001
$ws.Range("H106").Value = 7999
$ws.Range("J106").Value = 7999
$ws.Range("L106").Value = 23997
$ws.Range("N106").Value = -25889
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 21999.2
$ws.Range("J44").Value = 24999.334
$ws.Range("L44").Value = 24999.334
$ws.Range("N44").Value = -26191.334
$ws.Range("H80").Value = 5968.6665
$ws.Range("J80").Value = 5968.6665
$ws.Range("L80").Value = 5968.6665
$ws.Range("N80").Value = -7964.6665
$ws.Range("H83").Value = 5968.6665
$ws.Range("J83").Value = 5968.6665
$ws.Range("L83").Value = 29843.3325
$ws.Range("N83").Value = -39827.3325
$ws.Range("H126").Value = 1940.3334
$ws.Range("I126").Value = 1900.3572
$ws.Range("K126").Value = 5701.071599999999
$ws.Range("M126").Value = -3231.071599999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1592.5625
$ws.Range("I16").Value = 1574.8695
$ws.Range("K16").Value = 1574.8695
$ws.Range("M16").Value = -1404.8695
$ws.Range("H61").Value = 3113.2666
$ws.Range("I61").Value = 3113.2666
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3113.2666
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2911.2666
$ws.Range("N61").ClearContents()
$ws.Range("H82").Value = 4610.7827
$ws.Range("J82").Value = 4383.143
$ws.Range("L82").Value = 4383.143
$ws.Range("N82").Value = -5105.143
$ws.Range("H85").Value = 4610.7827
$ws.Range("J85").Value = 4383.143
$ws.Range("L85").Value = 4383.143
$ws.Range("N85").Value = -6879.143
$ws.Range("H113").Value = 3113.2666
$ws.Range("I113").Value = 3113.2666
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3113.2666
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -943.2665999999999
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 60666.332
$ws.Range("I132").Value = 4000
$ws.Range("K132").Value = 12000
$ws.Range("M132").Value = -9470
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1547
$ws.Range("I113").Value = 987.3333
$ws.Range("K113").Value = 2961.9999
$ws.Range("M113").Value = -791.9998999999998
$ws.Range("H116").Value = 45677
$ws.Range("J116").Value = 45677
$ws.Range("L116").Value = 45677
$ws.Range("N116").Value = -54855
